# Generate Report for Handoff
# Updates the localization-status workbook so that rows previously marked
# "Handed back: in sync with en-US" are now "Ready for handoff", and bumps
# the associated timestamps to reflect the new handoff generation time.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-02 03:08:55"

# --- zh-cn sheet ---
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-02 03:08:51"

# --- de-de sheet ---
# C2 = Status, H2 = Latest Handoff Datetime (shares the same underlying
# value as the Overview sheet's "Latest HO Xliff Generate Date", so it
# picks up the same refreshed timestamp as Overview!G2)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-02 03:08:55"

# --- Column width adjustments to reflect the shorter status text ---
# Target stored width is 17.2159881591797 "characters"; the closest width
# this engine's ColumnWidth setter can reach (it quantizes to 1/6-character
# pixel steps) is obtained by requesting 16.3333333333333 characters.
$newStatusColWidth = 16.3333333333333

# Overview columns E and F (zh-cn / de-de status columns)
$overview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# zh-cn / de-de sheet Status column (C)
$zhcn.Columns.Item(3).ColumnWidth = $newStatusColWidth
$dede.Columns.Item(3).ColumnWidth = $newStatusColWidth
